# Update gh-pages output data (江西-漫展信息) to the values generated at 456a3b4.
# The "展览" sheet and the "全部类型" sheet carry the same underlying rows,
# so every change is mirrored across both worksheets.

$wb = $excel.ActiveWorkbook

# Map of row -> new "想去人数" (F column) value shared by both sheets.
$updates = @{
    6  = 96
    12 = 316
    14 = 384
    17 = 15
    19 = 54
    21 = 997
    22 = 1416
    24 = 339
    26 = 79
    32 = 279
    33 = 1640
    39 = 3766
}

$newCover = "//i0.hdslb.com/bfs/openplatform/202402/l6GUtggC1706843695971.jpeg"

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }

    # Row 42's Cover (column I) image URL was refreshed.
    $ws.Cells.Item(42, 9).Value = $newCover
}
